# Scheduled runner update: refresh market-price-derived columns (H..N)
# on the Kujata_Profits leve-profit sheets, per the latest price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value2 = 820.8570999999999
$ws.Range("I54").Value2 = 820.8570999999999
$ws.Range("J54").Value2 = 0
$ws.Range("K54").Value2 = 820.8570999999999
$ws.Range("L54").Value2 = 0
$ws.Range("M54").Value2 = -334.8570999999999
$ws.Range("N54").ClearContents()

$ws.Range("H106").Value2 = 6044.963
$ws.Range("I106").Value2 = 6162.077
$ws.Range("K106").Value2 = 6162.077
$ws.Range("M106").Value2 = -5531.077

$ws.Range("H112").Value2 = 2544.2222
$ws.Range("J112").Value2 = 3071.1428
$ws.Range("L112").Value2 = 9213.428400000001
$ws.Range("N112").Value2 = -11429.4284

$ws.Range("H137").Value2 = 1473.7428
$ws.Range("I137").Value2 = 1127.9412
$ws.Range("K137").Value2 = 3383.8236
$ws.Range("M137").Value2 = -833.8235999999997

$ws.Range("H138").Value2 = 531900.3
$ws.Range("J138").Value2 = 622812.3
$ws.Range("L138").Value2 = 1868436.9
$ws.Range("N138").Value2 = -1878716.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 2954.31
$ws.Range("I32").Value2 = 2428.7883
$ws.Range("J32").Value2 = 5932.2666
$ws.Range("K32").Value2 = 2428.7883
$ws.Range("L32").Value2 = 5932.2666
$ws.Range("M32").Value2 = -2141.7883
$ws.Range("N32").Value2 = -6506.2666

$ws.Range("H45").Value2 = 1319.1111
$ws.Range("I45").Value2 = 1127.4286
$ws.Range("J45").Value2 = 1990
$ws.Range("K45").Value2 = 1127.4286
$ws.Range("L45").Value2 = 1990
$ws.Range("M45").Value2 = -750.4286
$ws.Range("N45").Value2 = -2744

$ws.Range("H63").Value2 = 18520210
$ws.Range("I63").Value2 = 1507.425
$ws.Range("J63").Value2 = 71430780
$ws.Range("K63").Value2 = 1507.425
$ws.Range("L63").Value2 = 71430780
$ws.Range("M63").Value2 = -821.425
$ws.Range("N63").Value2 = -71432152

$ws.Range("H66").Value2 = 18520210
$ws.Range("I66").Value2 = 1507.425
$ws.Range("J66").Value2 = 71430780
$ws.Range("K66").Value2 = 7537.125
$ws.Range("L66").Value2 = 357153900
$ws.Range("M66").Value2 = -4105.125
$ws.Range("N66").Value2 = -357160764

$ws.Range("H74").Value2 = 1562.3529
$ws.Range("I74").Value2 = 1166.1538
$ws.Range("J74").Value2 = 2850
$ws.Range("K74").Value2 = 1166.1538
$ws.Range("L74").Value2 = 2850
$ws.Range("M74").Value2 = -292.1538
$ws.Range("N74").Value2 = -4598

$ws.Range("H77").Value2 = 1562.3529
$ws.Range("I77").Value2 = 1166.1538
$ws.Range("J77").Value2 = 2850
$ws.Range("K77").Value2 = 5830.769
$ws.Range("L77").Value2 = 14250
$ws.Range("M77").Value2 = -1462.769
$ws.Range("N77").Value2 = -22986

$ws.Range("H92").Value2 = 16750
$ws.Range("J92").Value2 = 16750
$ws.Range("L92").Value2 = 16750
$ws.Range("N92").Value2 = -21742

$ws.Range("H122").Value2 = 2497.5
$ws.Range("I122").Value2 = 2176.5
$ws.Range("J122").Value2 = 3139.5
$ws.Range("K122").Value2 = 6529.5
$ws.Range("L122").Value2 = 9418.5
$ws.Range("M122").Value2 = -4079.5
$ws.Range("N122").Value2 = -14318.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value2 = 2641.4583
$ws.Range("I86").Value2 = 2470.5293
$ws.Range("J86").Value2 = 3056.5715
$ws.Range("K86").Value2 = 2470.5293
$ws.Range("L86").Value2 = 3056.5715
$ws.Range("M86").Value2 = -1347.5293
$ws.Range("N86").Value2 = -5302.5715

$ws.Range("H89").Value2 = 2641.4583
$ws.Range("I89").Value2 = 2470.5293
$ws.Range("J89").Value2 = 3056.5715
$ws.Range("K89").Value2 = 12352.6465
$ws.Range("L89").Value2 = 15282.8575
$ws.Range("M89").Value2 = -6736.646500000001
$ws.Range("N89").Value2 = -26514.8575

$ws.Range("H99").Value2 = 166667970
$ws.Range("I99").Value2 = 200001220
$ws.Range("J99").Value2 = 1700
$ws.Range("K99").Value2 = 200001220
$ws.Range("L99").Value2 = 1700
$ws.Range("M99").Value2 = -199999722
$ws.Range("N99").Value2 = -4696

$ws.Range("H105").Value2 = 50001530
$ws.Range("I105").Value2 = 58825010
$ws.Range("J105").Value2 = 1803.6666
$ws.Range("K105").Value2 = 58825010
$ws.Range("L105").Value2 = 1803.6666
$ws.Range("M105").Value2 = -58823263
$ws.Range("N105").Value2 = -5297.6666

$ws.Range("H119").Value2 = 45400
$ws.Range("J119").Value2 = 45400
$ws.Range("L119").Value2 = 45400
$ws.Range("N119").Value2 = -55076

$ws.Range("H140").Value2 = 23856.47
$ws.Range("J140").Value2 = 23856.47
$ws.Range("L140").Value2 = 23856.47
$ws.Range("N140").Value2 = -34216.47

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 45455708
$ws.Range("I16").Value2 = 66667810
$ws.Range("J16").Value2 = 1199.8572
$ws.Range("K16").Value2 = 66667810
$ws.Range("L16").Value2 = 1199.8572
$ws.Range("M16").Value2 = -66667523
$ws.Range("N16").Value2 = -1773.8572

$ws.Range("H31").Value2 = 1604.0294
$ws.Range("I31").Value2 = 1538.6562
$ws.Range("J31").Value2 = 2650
$ws.Range("K31").Value2 = 1538.6562
$ws.Range("L31").Value2 = 2650
$ws.Range("M31").Value2 = -1243.6562
$ws.Range("N31").Value2 = -3240

$ws.Range("H34").Value2 = 1604.0294
$ws.Range("I34").Value2 = 1538.6562
$ws.Range("J34").Value2 = 2650
$ws.Range("K34").Value2 = 1538.6562
$ws.Range("L34").Value2 = 2650
$ws.Range("M34").Value2 = -1336.6562
$ws.Range("N34").Value2 = -3054

$ws.Range("H105").Value2 = 899.8889
$ws.Range("J105").Value2 = 1055
$ws.Range("L105").Value2 = 1055
$ws.Range("N105").Value2 = -4549

$ws.Range("H113").Value2 = 45455708
$ws.Range("I113").Value2 = 66667810
$ws.Range("J113").Value2 = 1199.8572
$ws.Range("K113").Value2 = 66667810
$ws.Range("L113").Value2 = 1199.8572
$ws.Range("M113").Value2 = -66665640
$ws.Range("N113").Value2 = -5539.8572

$ws.Range("H122").Value2 = 1067.8
$ws.Range("I122").Value2 = 1356
$ws.Range("J122").Value2 = 875.6667
$ws.Range("K122").Value2 = 4068
$ws.Range("L122").Value2 = 2627.0001
$ws.Range("M122").Value2 = -1618
$ws.Range("N122").Value2 = -7527.0001

$ws.Range("H132").Value2 = 1223.4286
$ws.Range("I132").Value2 = 836.9167
$ws.Range("K132").Value2 = 2510.7501
$ws.Range("M132").Value2 = 19.2498999999998

$ws.Range("H134").Value2 = 31252610
$ws.Range("I134").Value2 = 2986.9
$ws.Range("K134").Value2 = 8960.700000000001
$ws.Range("M134").Value2 = -6425.700000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value2 = 75
$ws.Range("I8").Value2 = 75
$ws.Range("K8").Value2 = 225
$ws.Range("M8").Value2 = -86

$ws.Range("H41").Value2 = 2166.6667
$ws.Range("I41").Value2 = 0
$ws.Range("J41").Value2 = 2166.6667
$ws.Range("K41").Value2 = 0
$ws.Range("L41").Value2 = 6500.000100000001
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value2 = -7176.000100000001

$ws.Range("H131").Value2 = 12195981
$ws.Range("I131").Value2 = 55556040
$ws.Range("J131").Value2 = 964
$ws.Range("K131").Value2 = 166668120
$ws.Range("L131").Value2 = 2892
$ws.Range("M131").Value2 = -166663080
$ws.Range("N131").Value2 = -12972

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value2 = 5000
$ws.Range("J53").Value2 = 5000
$ws.Range("L53").Value2 = 5000
$ws.Range("N53").Value2 = -6262

$ws.Range("H102").Value2 = 1180.6666
$ws.Range("I102").Value2 = 1300.4445
$ws.Range("K102").Value2 = 1300.4445
$ws.Range("M102").Value2 = 321.5554999999999

$ws.Range("H134").Value2 = 23512.076
$ws.Range("J134").Value2 = 25054.75
$ws.Range("L134").Value2 = 75164.25
$ws.Range("N134").Value2 = -80234.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value2 = 0
$ws.Range("I17").Value2 = 0
$ws.Range("K17").Value2 = 0
$ws.Range("M17").ClearContents()

$ws.Range("H57").Value2 = 10000
$ws.Range("J57").Value2 = 10000
$ws.Range("L57").Value2 = 10000
$ws.Range("N57").Value2 = -11132

$ws.Range("H61").Value2 = 1115.1666
$ws.Range("I61").Value2 = 1058.2
$ws.Range("K61").Value2 = 1058.2
$ws.Range("M61").Value2 = -856.2

$ws.Range("H113").Value2 = 1115.1666
$ws.Range("I113").Value2 = 1058.2
$ws.Range("K113").Value2 = 1058.2
$ws.Range("M113").Value2 = 1111.8

$ws.Range("H122").Value2 = 31252208
$ws.Range("I122").Value2 = 50002132
$ws.Range("J122").Value2 = 2335
$ws.Range("K122").Value2 = 150006396
$ws.Range("L122").Value2 = 7005
$ws.Range("M122").Value2 = -150003946
$ws.Range("N122").Value2 = -11905

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value2 = 39649.75
$ws.Range("J109").Value2 = 36085.668
$ws.Range("L109").Value2 = 36085.668
$ws.Range("N109").Value2 = -38859.668

$ws.Range("H131").Value2 = 29306
$ws.Range("J131").Value2 = 29306
$ws.Range("L131").Value2 = 29306
$ws.Range("N131").Value2 = -39386

$ws.Range("H132").Value2 = 1969.1428
$ws.Range("I132").Value2 = 1630
$ws.Range("K132").Value2 = 4890
$ws.Range("M132").Value2 = -2360

$ws.Range("H136").Value2 = 1199.7142
$ws.Range("I136").Value2 = 1056.8
$ws.Range("K136").Value2 = 3170.4
$ws.Range("M136").Value2 = -620.3999999999996
